$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header "Modelo" in F1, matching the style used by the other headers
$ws.Range("F1").Value = "Modelo"
$ws.Range("A1").Copy()
$ws.Range("F1").PasteSpecial(-4122)
$ws.Range("F1").Value = "Modelo"
$excel.CutCopyMode = 0

$modelText = "MultiOutputRegressor(estimator=GridSearchCV(cv=5,`n                                            estimator=Pipeline(steps=[('model',`n                                                                       GradientBoostingRegressor())]),`n                                            param_grid={'model__max_depth': [3,`n                                                                             5,`n                                                                             7],`n                                                        'model__n_estimators': [50,`n                                                                                100,`n                                                                                150]},`n                                            scoring='neg_mean_squared_error'))"

$ws.Range("F2").Value = $modelText
$ws.Range("F3").Value = $modelText
$ws.Range("F4").Value = $modelText
$ws.Range("F5").Value = $modelText

# Update the slightly adjusted MSE/MAE values for rows 3-5
$ws.Range("B3").Value = 0.07514644587374572
$ws.Range("D3").Value = 0.2119198634755614

$ws.Range("B4").Value = 0.0421553411937144
$ws.Range("D4").Value = 0.1361288253571674

$ws.Range("B5").Value = 0.07796894984218621
$ws.Range("D5").Value = 0.1911874935925041
